# The scenario rate values (Data/scenarios.xlsx) were stored as raw decimal
# fractions (e.g. -0.0234...). The rate computation in the filters now
# expects values expressed on a 0-100 scale, so every rate cell is rescaled
# by a factor of 100. Only the data block (rows 2-8, columns C:AF) holds
# rates; column A is the scenario label and column B is blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rates = $ws.Range("C2:AF8")

foreach ($cell in $rates.Cells) {
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value = $current * 100
    }
}
